$d = $word.ActiveDocument

# --- 1. Remove the bookmark from its old location (it will be moved to the
#        new edit point, mirroring Word's automatic "_GoBack" tracking) and
#        normalize the run split left behind. ---
$pReset = $d.Paragraphs(91)
if ($pReset.Range.Text -notmatch "reset = 0;") {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -match "reset = 0;") {
            $pReset = $d.Paragraphs($i)
            break
        }
    }
}
$resetXml = '<?xml version="1.0" encoding="utf-8"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>reset = 0;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$pReset.Range.InsertXML($resetXml)

# --- 2. Type the new first line of text ("Myname is swanil") into what was
#        an empty paragraph, including the automatic spell-check proofErr
#        markers around "Myname", and re-seat "_GoBack" at this, the most
#        recent edit location. ---
$p1 = $d.Paragraphs(1)
$introXml = '<?xml version="1.0" encoding="utf-8"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Myname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is swanil</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p1.Range.InsertXML($introXml)
